# Draft mapping addition: adds a new "Mapping: Spécification métier vers
# l'extension ROR ContactConfidentialityLevel" column to the Elements sheet,
# fills in the mapping value for the Extension.value[x] element, and bumps
# the StructureDefinition's Date metadata.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: refresh the "Date" property -----------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-12T09:15:29+00:00"

# --- Elements sheet: add the new mapping column (AL) --------------------
$ws = $wb.Worksheets.Item("Elements")

# Copy header/data formatting from the last existing mapping column (AK)
# onto the new column (AL) so it matches the other "Mapping: ..." columns.
$ws.Range("AK1:AK6").Copy()
$ws.Range("AL1:AL6").PasteSpecial(-4122)

$ws.Range("AL1").Value = "Mapping: Spécification métier vers l'extension ROR ContactConfidentialityLevel"
$ws.Range("AL2").Value = ""
$ws.Range("AL3").Value = ""
$ws.Range("AL4").Value = ""
$ws.Range("AL5").Value = ""
$ws.Range("AL6").Value = "niveauConfidentialite"

# Column AL width to match the other wide mapping-text columns.
$ws.Columns.Item(38).ColumnWidth = 83.7734375
